$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (D column explanations + B28 answer) ---
$ws.Range("D4").Value = '유병기간의 추측에 중요한 질문이며 35세 여성이 대학생때부터 설사를 했으므로 만성 설사를 시사하는 질문입니다. '
$ws.Range("D5").Value = '설사가 악화된 기간, 설사시 시 동반 증상을 확인합니다.'
$ws.Range("D6").Value = '하루중 설사 횟수가 3회 이상으로 만성 설사에 합당합니다. '
$ws.Range("D7").Value = '만성 설사 중 설사양이 많을 때는 삼투성, 분비성 설사를 의심해야하며, 많지 않을 때는 과민성 대장증후군. 염증성 질환에 의한 설사를 감별해야햡니다. '
$ws.Range("D8").Value = '기질성 설사의 경우는 밤에 설사를 하는 경우도 흔하며 기능성 원인에 의한 설사일 경우 자다가 깨서 급하게 설사를 보는 경우는 드물기 때문에.  기능성 설사를 시사하는 소견입니다. '
$ws.Range("D9").Value = '설사가 식사와 관련성을 확인하며 식후에 30분이내 설사시 약제, 유당섭취력이 없다면 과민성 장증후군이 가장흔한 원인입니다. '
$ws.Range("D10").Value = '점맥변이나 혈변은 장에 염증에 의한 설사이나 상기 증상이 없으므로 대장에 염증에 의한 염증성 설사의 가능성은 낮습니다. '
$ws.Range("D11").Value = '긴장하면 설사를 하는 것은 과민성 대장 증후군의 가능성이 높습니다. '
$ws.Range("D12").Value = '긴장하거나 스트레스로 설사를 하는 것은 과민성 대장 증후군의 가능성이 높습니다. '
$ws.Range("D13").Value = '항문농양, 치루 등의 동반증상을 확인할 수 있으면 항문질환이 동반되면 크론병과 같은 염증성 장질환의 감별진단에 도움이 됩니다. 본 환자는 상기 증상이 없으므로  크론병의 가능성은 낮습니다. '
$ws.Range("D14").Value = '만성설사의 원인 중 발열, 식욕소실, 체중감소, 소변이상이 없으므로  염증성 장질환에 의한 설사는 가능성이 낮읍니다. '
$ws.Range("D15").Value = '증상은 오래되었지만 3년전 대장내시경 검사가 정상이었으므로 기능성 설사일  가능성이 높습니다. '
$ws.Range("D16").Value = 'Vagotome, gastrectomy, blind loop with bacterial overgrowth 수술력이 없으므로 앞 수술과 관련된  설사, 당뇨, 갑상성 기능항진증등의 의한 설사일 가능성은 없습니다.'
$ws.Range("D17").Value = '아스피린 등은 장기 복용시 설사를 유발할 수 있으나 본 환자는 투약력이 없으므로  약제유발성 설사의 가능성은 낮습니다. '
$ws.Range("D18").Value = '알코올은 설사를 유발할 가능성이 있습니다. '
$ws.Range("D19").Value = '흡연이 크론병과 같은 염증성 장질환의 발생을 증가시기는 하지만 끊은 상태이므로 염증성 장질환의 가눙성은 낮습니다. '
$ws.Range("D20").Value = '만성 설사의 원인 중  아메바나 지아르디아 기생충 감염은 동남아 여행 후 발생할 수 있으나 본 환자는 동남아  출장후 몸이 더 불편해지지는 않았으므로  아메바, 지아르디아 등의 기생충감염의 가능성은 없습니다. '
$ws.Range("D21").Value = '만성 염증성장질환등은 가족력이 있으므로 염증성 장질환에 의한 설사 가능성은 낮습니다. '
$ws.Range("D24").Value = '물론 검사를 해보아야 겠지만 과민성 장증군, 알코올 유발성 설사, 음식알러지 등을 감별해야 할 것 같습니다. 너무 걱정하지 마세요.'
$ws.Range("D25").Value = '혈액검사, 변검사, 대장내시경검사를 받게 되실 겁니다. '
$ws.Range("B28").Value = '만성 설사(4주 이상 설사)의 흔한 원인 5가지는  과민성 장증후군, 음식 알러지, 염증성 장질환, 감염성 설사, 유당불내성  등이 있다. '

# --- Row height adjustments (content grew, rows need more height) ---
$ws.Rows.Item(7).RowHeight = 39.75
$ws.Rows.Item(8).RowHeight = 37.5
$ws.Rows.Item(9).RowHeight = 57.75
$ws.Rows.Item(10).RowHeight = 48
$ws.Rows.Item(12).RowHeight = 41.25
$ws.Rows.Item(13).RowHeight = 40.5
$ws.Rows.Item(16).RowHeight = 30.75
$ws.Rows.Item(20).RowHeight = 37.5

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 107
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 108.71428571428571

# --- Selection / view update ---
$ws.Range("D12").Select()
